$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update A12: clarify that deleting is now implemented in the backend ---
$ws.Range("A12").Value = "Deleting: (implemented in backend - now just need front end to handle)"

# --- Reshuffle the TODO items that used to sit in rows 16-21 ---
# Row 21's old note ("Be able to unsubscribe from trucks") becomes the new
# A16 (re-worded/indented); a new clarifying note is inserted at A18; and the
# remaining rows 17-20 each shift down by two rows (-> 19-22), while the old
# A16 shifts down by one row (-> 17). Read each source before it is
# overwritten, working from the bottom up.
$ws.Range("A22").Value = $ws.Range("A20").Value2
$ws.Range("A21").Value = $ws.Range("A19").Value2
$ws.Range("A20").Value = $ws.Range("A18").Value2
$ws.Range("A19").Value = $ws.Range("A17").Value2
$ws.Range("A17").Value = $ws.Range("A16").Value2
$ws.Range("A16").Value = "    Be able to unsubscribe from trucks"
$ws.Range("A18").Value = "     (implemented in backend - now just need front end to handle)"

# Row 17 ("Implement userDetails...") now carries the yellow highlight that
# used to sit on row 16
$ws.Range("A17").Interior.Color = 65535

# Row 16 ("Be able to unsubscribe from trucks") no longer carries a highlight
$ws.Range("A16").Interior.Pattern = -4142

# Row 18 ("(implemented in backend...)") has no special fill
$ws.Range("A18").Interior.Pattern = -4142

# Select A16, matching the final cursor position left by the edit
$ws.Range("A16").Select()
